# Edit: add January entries (Jan 2-5, 2026) and update Total sheet summary figures.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) "January" sheet: insert 32 new homework-submission rows above the
#    existing 5 (they end up at the bottom, rows 34-38, unchanged).
# ---------------------------------------------------------------------------
$wsJan = $wb.Worksheets.Item("January")

$janRows = @(
    @{Row=2; A="https://lms.partaonline.ru/mentor/homeworks/650818"; B="ЕГЭ"; C=46027; D="23:08:00"; E="January"; F=4; G=40},
    @{Row=3; A="https://lms.partaonline.ru/mentor/homeworks/650650"; B="ЕГЭ"; C=46027; D="19:33:00"; E="January"; F=4; G=40},
    @{Row=4; A="https://lms.partaonline.ru/mentor/homeworks/650649"; B="ЕГЭ"; C=46027; D="19:32:00"; E="January"; F=3; G=30},
    @{Row=5; A="https://lms.partaonline.ru/mentor/homeworks/650648"; B="ЕГЭ"; C=46027; D="19:32:00"; E="January"; F=4; G=40},
    @{Row=6; A="https://lms.partaonline.ru/mentor/homeworks/650645"; B="ЕГЭ"; C=46027; D="19:31:00"; E="January"; F=7; G=70},
    @{Row=7; A="https://lms.partaonline.ru/mentor/homeworks/650644"; B="ЕГЭ"; C=46027; D="19:30:00"; E="January"; F=5; G=50},
    @{Row=8; A="https://lms.partaonline.ru/mentor/homeworks/650642"; B="ЕГЭ"; C=46027; D="19:30:00"; E="January"; F=4; G=40},
    @{Row=9; A="https://lms.partaonline.ru/mentor/homeworks/650640"; B="ЕГЭ"; C=46027; D="19:28:00"; E="January"; F=1; G=10},
    @{Row=10; A="https://lms.partaonline.ru/mentor/homeworks/650570"; B="ЕГЭ"; C=46027; D="18:22:00"; E="January"; F=4; G=40},
    @{Row=11; A="https://lms.partaonline.ru/mentor/homeworks/650496"; B="ОГЭ"; C=46027; D="16:05:00"; E="January"; F=1; G=18},
    @{Row=12; A="https://lms.partaonline.ru/mentor/homeworks/650490"; B="ЕГЭ"; C=46027; D="16:04:00"; E="January"; F=4; G=40},
    @{Row=13; A="https://lms.partaonline.ru/mentor/homeworks/650470"; B="ОГЭ"; C=46027; D="16:04:00"; E="January"; F=1; G=18},
    @{Row=14; A="https://lms.partaonline.ru/mentor/homeworks/650342"; B="ЕГЭ"; C=46027; D="12:15:00"; E="January"; F=7; G=70},
    @{Row=15; A="https://lms.partaonline.ru/mentor/homeworks/650219"; B="ЕГЭ"; C=46027; D="12:11:00"; E="January"; F=3; G=30},
    @{Row=16; A="https://lms.partaonline.ru/mentor/homeworks/650171"; B="ЕГЭ"; C=46027; D="12:10:00"; E="January"; F=5; G=50},
    @{Row=17; A="https://lms.partaonline.ru/mentor/homeworks/650150"; B="ЕГЭ"; C=46027; D="12:09:00"; E="January"; F=4; G=40},
    @{Row=18; A="https://lms.partaonline.ru/mentor/homeworks/650093"; B="ЕГЭ"; C=46027; D="12:04:00"; E="January"; F=6; G=60},
    @{Row=19; A="https://lms.partaonline.ru/mentor/homeworks/650079"; B="ЕГЭ"; C=46027; D="11:56:00"; E="January"; F=4; G=40},
    @{Row=20; A="https://lms.partaonline.ru/mentor/homeworks/650004"; B="ОГЭ"; C=46027; D="11:55:00"; E="January"; F=2; G=36},
    @{Row=21; A="https://lms.partaonline.ru/mentor/homeworks/649957"; B="ОГЭ"; C=46027; D="11:48:00"; E="January"; F=1; G=18},
    @{Row=22; A="https://lms.partaonline.ru/mentor/homeworks/649928"; B="ЕГЭ"; C=46027; D="11:48:00"; E="January"; F=5; G=50},
    @{Row=23; A="https://lms.partaonline.ru/mentor/homeworks/649807"; B="ОГЭ"; C=46026; D="13:29:00"; E="January"; F=3; G=54},
    @{Row=24; A="https://lms.partaonline.ru/mentor/homeworks/649792"; B="ЕГЭ"; C=46026; D="13:26:00"; E="January"; F=5; G=50},
    @{Row=25; A="https://lms.partaonline.ru/mentor/homeworks/649674"; B="ЕГЭ"; C=46026; D="13:24:00"; E="January"; F=7; G=70},
    @{Row=26; A="https://lms.partaonline.ru/mentor/homeworks/649593"; B="ЕГЭ"; C=46026; D="13:19:00"; E="January"; F=1; G=10},
    @{Row=27; A="https://lms.partaonline.ru/mentor/homeworks/649508"; B="ОГЭ"; C=46026; D="13:23:00"; E="January"; F=1; G=18},
    @{Row=28; A="https://lms.partaonline.ru/mentor/homeworks/649480"; B="ОГЭ"; C=46026; D="13:22:00"; E="January"; F=3; G=54},
    @{Row=29; A="https://lms.partaonline.ru/mentor/homeworks/649479"; B="ЕГЭ"; C=46026; D="13:19:00"; E="January"; F=1; G=10},
    @{Row=30; A="https://lms.partaonline.ru/mentor/homeworks/649450"; B="ЕГЭ"; C=46026; D="13:21:00"; E="January"; F=5; G=50},
    @{Row=31; A="https://lms.partaonline.ru/mentor/homeworks/649392"; B="ЕГЭ"; C=46026; D="13:20:00"; E="January"; F=4; G=40},
    @{Row=32; A="https://lms.partaonline.ru/mentor/homeworks/649367"; B="ЕГЭ"; C=46026; D="13:19:00"; E="January"; F=4; G=40},
    @{Row=33; A="https://lms.partaonline.ru/mentor/homeworks/649342"; B="ЕГЭ"; C=46026; D="13:18:00"; E="January"; F=10; G=100},
    @{Row=34; A="https://lms.partaonline.ru/mentor/homeworks/649239"; B="ЕГЭ"; C=46025; D="12:38:00"; E="January"; F=4; G=40},
    @{Row=35; A="https://lms.partaonline.ru/mentor/homeworks/649065"; B="ЕГЭ"; C=46025; D="12:38:00"; E="January"; F=1; G=10},
    @{Row=36; A="https://lms.partaonline.ru/mentor/homeworks/649002"; B="ЕГЭ"; C=46025; D="12:37:00"; E="January"; F=10; G=100},
    @{Row=37; A="https://lms.partaonline.ru/mentor/homeworks/648780"; B="ЕГЭ"; C=46024; D="12:38:00"; E="January"; F=4; G=40},
    @{Row=38; A="https://lms.partaonline.ru/mentor/homeworks/648723"; B="ОГЭ"; C=46024; D="12:39:00"; E="January"; F=1; G=18}
)

# Give every row in the new range C2:C38 the same date number-format as the
# existing data (copy format only, values are set explicitly below) so no
# stray style gets introduced.
$wsJan.Range("C2").Copy()
$wsJan.Range("C2:C38").PasteSpecial(-4122)

foreach ($r in $janRows) {
    $row = $r.Row
    $wsJan.Cells.Item($row, 1).Value = $r.A
    $wsJan.Cells.Item($row, 2).Value = $r.B
    $wsJan.Cells.Item($row, 3).Value = $r.C
    $wsJan.Cells.Item($row, 4).Value = $r.D
    $wsJan.Cells.Item($row, 5).Value = $r.E
    $wsJan.Cells.Item($row, 6).Value = $r.F
    $wsJan.Cells.Item($row, 7).Value = $r.G
}

# ---------------------------------------------------------------------------
# 2) "Total" sheet: refresh the aggregated January counters/money figures.
# ---------------------------------------------------------------------------
$wsTotal = $wb.Worksheets.Item("Total")

$wsTotal.Range("E4").Value = 82
$wsTotal.Range("E5").Value = 820
$wsTotal.Range("E6").Value = 8
$wsTotal.Range("E7").Value = 120
$wsTotal.Range("F10").Value = 143
$wsTotal.Range("F11").Value = 1534
$wsTotal.Range("E12").Value = 12090
$wsTotal.Range("F12").Value = 1534
$wsTotal.Range("E13").Value = 10518.3
$wsTotal.Range("F13").Value = 1334.58

# View-state tweaks captured on the "Total" sheet: zoomed to 145% with a new
# active selection.
$wsTotal.Activate()
$excel.ActiveWindow.Zoom = 145
$wsTotal.Range("I16").Select()
